# Auto update Excel log 2026-02-04 14:22:32
# Appends new sensor log rows captured by the monitoring pipeline
# to the PIR, Humidity and Temperature sheets.

$wb = $excel.ActiveWorkbook

$pirData = @(
    ,@("2026-02-04","14:21:27","14:00","Bathroom","No Motion","Inactive")
    ,@("2026-02-04","14:21:28","14:00","Bathroom","Motion Detected","Active")
    ,@("2026-02-04","14:21:35","14:00","Bathroom","No Motion","Inactive")
    ,@("2026-02-04","14:21:40","14:00","Bathroom","No Motion","Inactive")
    ,@("2026-02-04","14:21:43","14:00","Bathroom","Motion Detected","Active")
    ,@("2026-02-04","14:21:51","14:00","Bathroom","No Motion","Inactive")
    ,@("2026-02-04","14:21:56","14:00","Bathroom","No Motion","Inactive")
    ,@("2026-02-04","14:22:01","14:00","Bathroom","No Motion","Inactive")
    ,@("2026-02-04","14:22:06","14:00","Bathroom","No Motion","Inactive")
    ,@("2026-02-04","14:22:11","14:00","Bathroom","No Motion","Inactive")
    ,@("2026-02-04","14:22:16","14:00","Bathroom","No Motion","Inactive")
    ,@("2026-02-04","14:22:20","14:00","Bathroom","Motion Detected","Active")
)

$humidityData = @(
    ,@("2026-02-04","14:21:29","14:00","Bathroom","78.2%","Active")
    ,@("2026-02-04","14:21:34","14:00","Bathroom","78.4%","Active")
    ,@("2026-02-04","14:21:39","14:00","Bathroom","78.3%","Active")
    ,@("2026-02-04","14:21:44","14:00","Bathroom","78.1%","Active")
    ,@("2026-02-04","14:21:49","14:00","Bathroom","77.2%","Active")
    ,@("2026-02-04","14:21:54","14:00","Bathroom","78.2%","Active")
    ,@("2026-02-04","14:22:04","14:00","Bathroom","78.2%","Active")
    ,@("2026-02-04","14:22:09","14:00","Bathroom","77.3%","Active")
    ,@("2026-02-04","14:22:19","14:00","Bathroom","77.4%","Active")
    ,@("2026-02-04","14:22:24","14:00","Bathroom","78.3%","Active")
)

$temperatureData = @(
    ,@("2026-02-04","14:21:29","14:00","Bathroom","24.6C","Active")
    ,@("2026-02-04","14:21:34","14:00","Bathroom","24.7C","Active")
    ,@("2026-02-04","14:21:40","14:00","Bathroom","24.7C","Active")
    ,@("2026-02-04","14:21:44","14:00","Bathroom","24.7C","Active")
    ,@("2026-02-04","14:21:50","14:00","Bathroom","24.7C","Active")
    ,@("2026-02-04","14:21:55","14:00","Bathroom","24.7C","Active")
    ,@("2026-02-04","14:22:05","14:00","Bathroom","24.7C","Active")
    ,@("2026-02-04","14:22:10","14:00","Bathroom","24.6C","Active")
    ,@("2026-02-04","14:22:20","14:00","Bathroom","24.7C","Active")
    ,@("2026-02-04","14:22:25","14:00","Bathroom","24.7C","Active")
)


$wsPIR = $wb.Worksheets.Item("PIR")
$lastRow = 1
while ($wsPIR.Cells.Item($lastRow + 1, 1).Value() -ne $null -and $wsPIR.Cells.Item($lastRow + 1, 1).Value() -ne "") {
    $lastRow = $lastRow + 1
}
$r = $lastRow + 1
$textCols = @(1,2,3)
foreach ($row in $pirData) {
    for ($col = 1; $col -le 6; $col++) {
        $cell = $wsPIR.Cells.Item($r, $col)
        if ($textCols -contains $col) {
            $cell.NumberFormat = "@"
        }
        $cell.Value = $row[$col - 1]
    }
    $r = $r + 1
}


$wsHumidity = $wb.Worksheets.Item("Humidity")
$lastRow = 1
while ($wsHumidity.Cells.Item($lastRow + 1, 1).Value() -ne $null -and $wsHumidity.Cells.Item($lastRow + 1, 1).Value() -ne "") {
    $lastRow = $lastRow + 1
}
$r = $lastRow + 1
$textCols = @(1,2,3,5)
foreach ($row in $humidityData) {
    for ($col = 1; $col -le 6; $col++) {
        $cell = $wsHumidity.Cells.Item($r, $col)
        if ($textCols -contains $col) {
            $cell.NumberFormat = "@"
        }
        $cell.Value = $row[$col - 1]
    }
    $r = $r + 1
}


$wsTemperature = $wb.Worksheets.Item("Temperature")
$lastRow = 1
while ($wsTemperature.Cells.Item($lastRow + 1, 1).Value() -ne $null -and $wsTemperature.Cells.Item($lastRow + 1, 1).Value() -ne "") {
    $lastRow = $lastRow + 1
}
$r = $lastRow + 1
$textCols = @(1,2,3)
foreach ($row in $temperatureData) {
    for ($col = 1; $col -le 6; $col++) {
        $cell = $wsTemperature.Cells.Item($r, $col)
        if ($textCols -contains $col) {
            $cell.NumberFormat = "@"
        }
        $cell.Value = $row[$col - 1]
    }
    $r = $r + 1
}

